$wb = $excel.ActiveWorkbook

# --- Sheet "Estadisticos 1P" ---
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D2").Value = 0
$ws1.Range("E2").Value = 8
$ws1.Range("H2").Value = 7.5

$ws1.Range("D4").Value = 0
$ws1.Range("E4").Value = 5
$ws1.Range("H4").Value = 7.4

# --- Sheet "Estadisticos 2P" ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("D2").Value = 5
$ws2.Range("E2").Value = 11
$ws2.Range("F2").Value = 28
$ws2.Range("G2").Value = 71.79
$ws2.Range("H2").Value = 7.6

$ws2.Range("D3").Value = 4
$ws2.Range("E3").Value = 12
$ws2.Range("F3").Value = 25
$ws2.Range("G3").Value = 67.57
$ws2.Range("H3").Value = 6.5

$ws2.Range("D4").Value = 5
$ws2.Range("E4").Value = 5
$ws2.Range("F4").Value = 31
$ws2.Range("G4").Value = 86.11
$ws2.Range("H4").Value = 7.7

$ws2.Range("D5").Value = 6
$ws2.Range("E5").Value = 9
$ws2.Range("F5").Value = 26
$ws2.Range("G5").Value = 74.29
$ws2.Range("H5").Value = 7.2

# --- Sheet "Estadisticos Final" ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("D2").Value = 0
$ws3.Range("F2").Value = 32
$ws3.Range("G2").Value = 82.05
$ws3.Range("H2").Value = 7.7

$ws3.Range("H3").Value = 6.8

$ws3.Range("D4").Value = 0
$ws3.Range("E4").Value = 3
$ws3.Range("F4").Value = 33
$ws3.Range("G4").Value = 91.67
$ws3.Range("H4").Value = 8.1

$ws3.Range("H5").Value = 7.5

# --- Sheet "Rescatables" ---
$ws4 = $wb.Worksheets.Item("Rescatables")
$ws4.Range("A2").Value = 18330051920306
$ws4.Range("A3").Value = 18330051920429

$ws4.Range("B2").Value = "TEMOXTLE"
$ws4.Range("B3").Value = "GUERRA"

$ws4.Range("C2").Value = "LARA"
$ws4.Range("C3").Value = "OLMEDO"

$ws4.Range("D2").Value = "MADAI"
$ws4.Range("D3").Value = "PAOLA BETSABET"

$ws4.Range("E2").Value = "TEMAS DE BIOLOGÍA CONTEMPORÁNEA"
$ws4.Range("E3").Value = "TEMAS DE BIOLOGÍA CONTEMPORÁNEA"

$ws4.Range("F2").Value = "6ALCM"
$ws4.Range("F3").Value = "6BLCM"

$ws4.Range("G2").Value = 2
$ws4.Range("G3").Value = 2
